$wb = $excel.ActiveWorkbook

# Delete sheets that are no longer needed.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$wb.Worksheets.Item("Excercise 1").Delete()
$wb.Worksheets.Item("Exercise 2").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining "solutions" sheets.
$wb.Worksheets.Item("Shirt sales (i) (solutions)").Name = "Shirt sales (i)"
$wb.Worksheets.Item("Shirt sales (ii) (solutions)").Name = "Shirt sales (ii)"
